# Update gh-pages to output generated at 456a3b4
# Applies the F-column ("想去人数") numeric refreshes across all four sheets,
# and inserts a new event row into the "演出" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "藤川千爱2024演唱会" row into 演出 at row 41, pushing the
#    two rows that used to be 41/42 down to 42/43.
# ---------------------------------------------------------------------------
$wsPerf = $wb.Worksheets.Item("演出")
$wsPerf.Rows.Item(41).Insert()

$newRow = $wsPerf.Cells.Item(41, 1)
$newRow.Font.Bold = $true
$newRow.HorizontalAlignment = -4108
$newRow.VerticalAlignment = -4160
$newRow.Borders.LineStyle = 1
$newRow.Value = 40

$cB = $wsPerf.Cells.Item(41, 2)
$cB.NumberFormat = "@"
$cB.Value = "2024-07-21"
$cB.Style = "Normal"

$wsPerf.Cells.Item(41, 3).Value = "上海·藤川千爱2024演唱会"
$wsPerf.Cells.Item(41, 4).Value = "上海市普陀区宜昌路179号 万代南梦宫上海文化中心"

$cE = $wsPerf.Cells.Item(41, 5)
$cE.NumberFormat = "@"
$cE.Value = "2024.07.21 16:00-07.21 19:30"
$cE.Style = "Normal"

$wsPerf.Cells.Item(41, 6).Value = 1
$wsPerf.Cells.Item(41, 7).Value = 380
$wsPerf.Cells.Item(41, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85360"
$wsPerf.Cells.Item(41, 9).Value = "//i0.hdslb.com/bfs/openplatform/202405/Qhk2XOza1715248510067.jpeg"

# ---------------------------------------------------------------------------
# 2) Refresh the "想去人数" (F column) counters on every sheet.
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2 = 1420; 5 = 6675; 6 = 521; 8 = 36; 9 = 4553; 10 = 6774;
    12 = 220; 13 = 1382; 14 = 793; 15 = 111; 24 = 1046; 25 = 536;
    26 = 37; 27 = 27; 28 = 117; 30 = 1163; 31 = 28; 32 = 91; 34 = 3;
    37 = 514; 38 = 356; 39 = 33; 41 = 307; 42 = 1180; 43 = 518;
    44 = 60; 45 = 106
}
foreach ($r in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($r, 6).Value = $expoUpdates[$r]
}

$perfUpdates = @{
    2 = 12; 3 = 12; 7 = 516; 22 = 190; 31 = 750; 33 = 580; 35 = 89
}
foreach ($r in $perfUpdates.Keys) {
    $wsPerf.Cells.Item($r, 6).Value = $perfUpdates[$r]
}

$wsLocal = $wb.Worksheets.Item("本地生活")
$localUpdates = @{
    2 = 108; 4 = 710; 5 = 840; 6 = 605; 7 = 284; 8 = 1250; 9 = 1623
}
foreach ($r in $localUpdates.Keys) {
    $wsLocal.Cells.Item($r, 6).Value = $localUpdates[$r]
}

$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    2 = 12; 3 = 710; 4 = 1420; 6 = 840; 9 = 605; 10 = 605; 11 = 516;
    12 = 6675; 13 = 521; 15 = 36; 16 = 4553; 18 = 6774; 19 = 220;
    20 = 1382; 22 = 793; 23 = 111; 24 = 1250; 25 = 190; 29 = 1046;
    31 = 536; 32 = 37; 33 = 27; 34 = 117; 35 = 1163; 36 = 91;
    40 = 514; 41 = 580; 42 = 356; 43 = 33; 44 = 89; 45 = 307;
    46 = 518; 49 = 106
}
foreach ($r in $allUpdates.Keys) {
    $wsAll.Cells.Item($r, 6).Value = $allUpdates[$r]
}

Write-Host "edit complete"
